# PlanDeConvergenceIntegration.xlsx - "4 voltages monitoring added"
# Adds 13 new rows (104-116) to the "Liste_de_tâches" table on sheet 1,
# tweaks two existing rows (49 and 90), resizes the table / autofilter /
# conditional formatting / data validations accordingly, and updates the
# active selection on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------
# 1) Create the 13 new rows by copying the formatting of the last
#    existing data row (row 103) and the "blank" row 100 template so
#    that number formats / alignment match the existing table rows.
# ---------------------------------------------------------------------
$srcFull  = $ws.Range("B103:I103")   # fully populated template (state + % set)
$srcBlank = $ws.Range("B100:I100")   # blank État / % template

for ($r = 104; $r -le 116; $r++) {
    $srcFull.Copy($ws.Range("B$r`:I$r"))
}

# Rows 105, 106 and 114 have no État / % achevé value -> reuse the blank template
$srcBlank.Copy($ws.Range("B105:I105"))
$srcBlank.Copy($ws.Range("B106:I106"))
$srcBlank.Copy($ws.Range("B114:I114"))

$formula = '=IF(AND(Liste_de_tâches[[#This Row],[État ]]="Terminée",Liste_de_tâches[[#This Row],[% achevé]]=1),1,IF(ISBLANK(Liste_de_tâches[[#This Row],[Échéance ]]),2,IF(AND(Liste_de_tâches[[#This Row],[État ]]<>"Terminée",TODAY()>Liste_de_tâches[[#This Row],[Échéance ]]),3,2)))'
for ($r = 104; $r -le 116; $r++) {
    $ws.Range("H$r").Formula = $formula
}

# ---------------------------------------------------------------------
# 2) Fill in the task text (column B), following the exact order the
#    notes were typed in originally (drives shared-string ordering).
# ---------------------------------------------------------------------
$ws.Range("B104").Value = "Tester debuger le ping front "
$ws.Range("B105").Value = "Remplacer le servo moteur pour eviter tremblements"
$ws.Range("B106").Value = "Faire releve des pings en approche de portes"
$ws.Range("B107").Value = "Distinguer echo 0 de echo max"
$ws.Range("B108").Value = "surveiller le non decalage position encodeur VS roue"
$ws.Range("I108").Value = "via marque noire"
$ws.Range("B109").Value = "Monitorer l alim des moteurs (avant et apres regulateur)"
$ws.Range("B110").Value = "Stocker en BD les tensions "
$ws.Range("B111").Value = "Verifier / modifier le wait apres atteinte seuil encodeur pour prendre en compte la fin de mouvement "
$ws.Range("B112").Value = "Developper une fonction octave graph compare de 2 trajectoires avec heading"
$ws.Range("I112").Value = ' ApShowComparedStep(apRobot,det,actualPositions,"blue cyan: determined -- black red: actual")'
$ws.Range("B113").Value = "Java ne pas enregistrer les records scan recus en double"
$ws.Range("I111").Value = "ajout d un timer dans la loop en lieu et place du delay"
$ws.Range("I107").Value = "pour identifier les pbs d alimentation via une non montee du signal - a traiter via librairie NewPing"
$ws.Range("B114").Value = "Constituer un dossier de maintenance (procedure outils checklist,,,)"
$ws.Range("I104").Value = "semble irregulier VS ping back >> apres analyse SRF-05 return 0 aux environs de 3m (la doc parle de 4m) Front et Back sont similaires"
$ws.Range("B115").Value = "Developper une fonction test echo"
$ws.Range("I115").Value = "LoopPingFB,m + MySql repondent au besoin"
$ws.Range("B116").Value = "Etendre la carto au couloir et 1ere partie du salon "
$ws.Range("I116").Value = "a faire les jeux de tests"
$ws.Range("I113").Value = "to be checked"
$ws.Range("I109").Value = "a valider en reel"

# ---------------------------------------------------------------------
# 3) État (column D) and % achevé (column G) for every new row.
# ---------------------------------------------------------------------
$ws.Range("D104").Value = "En cours de réalisation"
$ws.Range("G104").Value = 0.75

$ws.Range("D107").Value = "Non commencée"

$ws.Range("D108").Value = "En cours de réalisation"
$ws.Range("G108").Value = 0.25

$ws.Range("D109").Value = "En cours de réalisation"
$ws.Range("G109").Value = 0.75

$ws.Range("D110").Value = "Terminée"
$ws.Range("G110").Value = 1

$ws.Range("D111").Value = "En cours de réalisation"
$ws.Range("G111").Value = 0.75

$ws.Range("D112").Value = "Terminée"
$ws.Range("G112").Value = 1

$ws.Range("D113").Value = "En cours de réalisation"
$ws.Range("G113").Value = 0.75

$ws.Range("D115").Value = "En cours de réalisation"
$ws.Range("G115").Value = 0.75

$ws.Range("D116").Value = "En cours de réalisation"
$ws.Range("G116").Value = 0.5

# ---------------------------------------------------------------------
# 4) Resize the table / autofilter to cover the new rows.
# ---------------------------------------------------------------------
$lo.Resize($ws.Range("B4:I116"))

# ---------------------------------------------------------------------
# 5) Update two pre-existing rows.
# ---------------------------------------------------------------------
$ws.Range("D49").Value = "Terminée"
$ws.Range("G49").Value = 1

$ws.Range("D90").Value = "En cours de réalisation"
$ws.Range("G90").Value = 0.5

# ---------------------------------------------------------------------
# 6) Refresh the active selection / scroll position on the sheet.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 100
$win.ScrollColumn = 3
$ws.Range("I107").Select()

$wb.Save()
